# Generate Report for Handoff
# Updates the "b.md" entries across the Overview, zh-cn and de-de sheets to
# reflect that the file is now "Ready for handoff" instead of "Handed back:
# in sync with en-US", and records the new handoff file / datetime.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-26-19 06:26:55"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-19 06:26:52"

foreach ($hl in $zhcn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$D$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# de-de sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-19 06:26:55"

foreach ($hl in $dede.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$D$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
